# Insert a new data row before row 312 (pushing existing rows 312..363 down to 313..364)
# and populate the new row with the new weekly record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(312).Insert()

$ws.Cells.Item(312, 1).Value = 6
$ws.Cells.Item(312, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(312, 3).Value = "Metropolitana"
$ws.Cells.Item(312, 4).Value = 45034
$ws.Cells.Item(312, 5).Value = 13
$ws.Cells.Item(312, 6).Value = 100112026
$ws.Cells.Item(312, 7).Value = "Haba"
$ws.Cells.Item(312, 8).Value = "Sin especificar"
$ws.Cells.Item(312, 9).Value = "Primera"
$ws.Cells.Item(312, 10).Value = 320
$ws.Cells.Item(312, 11).Value = 17000
$ws.Cells.Item(312, 12).Value = 18000
$ws.Cells.Item(312, 13).Value = 17625
$ws.Cells.Item(312, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(312, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(312, 16).Value = 705
$ws.Cells.Item(312, 17).Value = 25
$ws.Cells.Item(312, 18).Value = "Hortaliza"
